# Actualización automática 2025-09-08 09:31:10
# Updates sales figures on "VENTAS POR GRUPO" (by product group) and the
# matching monthly figures / budgets on "VENTA MENSUAL" for the advisor
# CASTRO ALCIVAR EDA MARIA.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws.Range("L4").Value  = 142.56
$ws.Range("E29").Value = 65.79000000000001
$ws.Range("D38").Value = 86.5
$ws.Range("L38").Value = 1309.61
$ws.Range("M38").Value = 3595.1
$ws.Range("L46").Value = 358.23
$ws.Range("M46").Value = -969.92
$ws.Range("L55").Value = 565.25
$ws.Range("M55").Value = 619.49

# Row 57 totals: "<count> de 55" tallies per column
$ws.Range("D57").Value = "1 de 55"
$ws.Range("E57").Value = "1 de 55"
$ws.Range("L57").Value = "4 de 55"

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$ws.Range("F4").Value  = 142.56
$ws.Range("G5").Value  = 3000
$ws.Range("G6").Value  = 4000
$ws.Range("G8").Value  = 3000
$ws.Range("G11").Value = 3000
$ws.Range("G12").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("G22").Value = 3000
$ws.Range("G24").Value = 1200
$ws.Range("G25").Value = 1500
$ws.Range("G26").Value = 3000
$ws.Range("G27").Value = 0
$ws.Range("G28").Value = 6327.1
$ws.Range("F29").Value = 65.79000000000001
$ws.Range("G29").Value = 6000
$ws.Range("G31").Value = 5000
$ws.Range("G32").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("G37").Value = 2300
$ws.Range("F38").Value = 4991.21
$ws.Range("G38").Value = 9679.110000000001
$ws.Range("G39").Value = 2000
$ws.Range("G41").Value = 3000
$ws.Range("G43").Value = 0
$ws.Range("G44").Value = 5600
$ws.Range("G45").Value = 4600
$ws.Range("F46").Value = -611.6900000000001
$ws.Range("G46").Value = 7468.67
$ws.Range("G47").Value = 2200
$ws.Range("G49").Value = 2000
$ws.Range("G50").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("G52").Value = 3000
$ws.Range("G53").Value = 2400
$ws.Range("G54").Value = 0
$ws.Range("F55").Value = 1184.74
$ws.Range("G55").Value = 2000

# Totals row
$ws.Range("F57").Value = 7099.71
$ws.Range("G57").Value = 85274.87999999999
